$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 38, pushing existing rows 38-51 down to 39-52.
$ws.Rows.Item(38).Insert()

# Populate the new row 38 with the weekly price-report entry that was added.
$ws.Cells.Item(38, 1).Value = 4
$ws.Cells.Item(38, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(38, 3).Value = "Los Lagos"
$ws.Cells.Item(38, 4).Value = 44960
$ws.Cells.Item(38, 5).Value = 10
$ws.Cells.Item(38, 6).Value = "Fruta"
$ws.Cells.Item(38, 7).Value = 100101
$ws.Cells.Item(38, 8).Value = "Berries"
$ws.Cells.Item(38, 9).Value = 100101001
$ws.Cells.Item(38, 10).Value = "Arándano (blue)"
$ws.Cells.Item(38, 11).Value = "Sin especificar"
$ws.Cells.Item(38, 12).Value = "Primera"
$ws.Cells.Item(38, 13).Value = 400
$ws.Cells.Item(38, 14).Value = 2000
$ws.Cells.Item(38, 15).Value = 2200
$ws.Cells.Item(38, 16).Value = 2100
$ws.Cells.Item(38, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(38, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(38, 19).Value = 1050
$ws.Cells.Item(38, 20).Value = 2
